$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '58.957.77'
$ws.Cells.Item(2, 5).Value = '  +2.52%  '

$ws.Cells.Item(3, 4).Value = '2.989.76'
$ws.Cells.Item(3, 5).Value = '  +1.37%  '

$ws.Cells.Item(4, 5).Value = '  +0.08%  '

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = '560.35'
$ws.Cells.Item(5, 4).ClearFormats()
$ws.Cells.Item(5, 5).Value = '  +0.90%  '

$ws.Cells.Item(6, 4).NumberFormat = "@"
$ws.Cells.Item(6, 4).Value = '136.76'
$ws.Cells.Item(6, 4).ClearFormats()
$ws.Cells.Item(6, 5).Value = '  +3.33%  '

$ws.Cells.Item(7, 5).Value = '  -0.06%  '

$ws.Cells.Item(8, 5).Value = '  +1.02%  '

$ws.Cells.Item(9, 4).Value = '2.977.49'
$ws.Cells.Item(9, 5).Value = '  +1.19%  '

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = '0.131'
$ws.Cells.Item(10, 4).ClearFormats()
$ws.Cells.Item(10, 5).Value = '  +2.55%  '

$ws.Cells.Item(11, 5).Value = '  +6.54%  '

$ws.Cells.Item(12, 5).Value = '  +1.78%  '

$ws.Cells.Item(13, 5).Value = '  +3.23%  '

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = '33.47'
$ws.Cells.Item(14, 4).ClearFormats()
$ws.Cells.Item(14, 5).Value = '  +2.02%  '

$ws.Cells.Item(15, 5).Value = '  +2.22%  '

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = '7.31'
$ws.Cells.Item(16, 4).ClearFormats()
$ws.Cells.Item(16, 5).Value = '  +7.25%  '

$ws.Cells.Item(17, 4).Value = '3.483.89'
$ws.Cells.Item(17, 5).Value = '  +1.46%  '

$ws.Cells.Item(18, 4).Value = '2.986.96'
$ws.Cells.Item(18, 5).Value = '  +1.61%  '

$ws.Cells.Item(19, 4).Value = '59.013.41'
$ws.Cells.Item(19, 5).Value = '  +2.61%  '

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = '426.01'
$ws.Cells.Item(20, 4).ClearFormats()
$ws.Cells.Item(20, 5).Value = '  +1.95%  '

$ws.Cells.Item(21, 4).NumberFormat = "@"
$ws.Cells.Item(21, 4).Value = '13.63'
$ws.Cells.Item(21, 4).ClearFormats()
$ws.Cells.Item(21, 5).Value = '  +3.70%  '

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = '0.722'
$ws.Cells.Item(22, 4).ClearFormats()
$ws.Cells.Item(22, 5).Value = '  +5.70%  '

$ws.Cells.Item(23, 5).Value = '  +1.76%  '

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = '13.25'
$ws.Cells.Item(24, 4).ClearFormats()
$ws.Cells.Item(24, 5).Value = '  +1.82%  '

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = '80.33'
$ws.Cells.Item(25, 4).ClearFormats()
$ws.Cells.Item(25, 5).Value = '  +1.35%  '

$ws.Cells.Item(26, 5).Value = '  -0.04%  '

$ws.Cells.Item(27, 5).Value = '  +0.15%  '

$ws.Cells.Item(28, 5).Value = '  +9.30%  '

$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = '2.53'
$ws.Cells.Item(29, 4).ClearFormats()
$ws.Cells.Item(29, 5).Value = '  +1.83%  '

$ws.Cells.Item(30, 5).Value = '  +2.72%  '

$ws.Cells.Item(31, 5).Value = '  +1.80%  '

$ws.Cells.Item(32, 5).Value = '  -1.06%  '

$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = '0.0990'
$ws.Cells.Item(33, 4).ClearFormats()
$ws.Cells.Item(33, 5).Value = '  -2.70%  '

$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = '0.993'
$ws.Cells.Item(34, 4).ClearFormats()
$ws.Cells.Item(34, 5).Value = '  +6.04%  '

$ws.Cells.Item(35, 5).Value = '  +4.96%  '

$ws.Cells.Item(36, 4).Value = '0.0₃0751'
$ws.Cells.Item(36, 5).Value = '  +9.54%  '

$ws.Cells.Item(37, 5).Value = '  -0.90%  '

$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = '48.70'
$ws.Cells.Item(38, 4).ClearFormats()
$ws.Cells.Item(38, 5).Value = '  +0.35%  '

$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = '8.68'
$ws.Cells.Item(39, 4).ClearFormats()
$ws.Cells.Item(39, 5).Value = '  +2.14%  '

$ws.Cells.Item(40, 5).Value = '  +5.65%  '

$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = '396.91'
$ws.Cells.Item(41, 4).ClearFormats()
$ws.Cells.Item(41, 5).Value = '  +4.62%  '

$ws.Cells.Item(42, 5).Value = '  +0.38%  '

$ws.Cells.Item(43, 4).Value = '2.746.91'
$ws.Cells.Item(43, 5).Value = '  +3.47%  '

$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = '0.107'
$ws.Cells.Item(44, 4).ClearFormats()
$ws.Cells.Item(44, 5).Value = '  -1.13%  '

$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = '0.250'
$ws.Cells.Item(45, 4).ClearFormats()
$ws.Cells.Item(45, 5).Value = '  +3.97%  '

$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = '35.01'
$ws.Cells.Item(47, 4).ClearFormats()
$ws.Cells.Item(47, 5).Value = '  +23.00%  '

$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = '122.54'
$ws.Cells.Item(48, 4).ClearFormats()
$ws.Cells.Item(48, 5).Value = '  -0.17%  '

$ws.Cells.Item(49, 5).Value = '  +0.60%  '

$ws.Cells.Item(50, 5).Value = '  +0.12%  '

$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = '23.24'
$ws.Cells.Item(51, 4).ClearFormats()
$ws.Cells.Item(51, 5).Value = '  -0.76%  '
